$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
# Row 17 (Leve Item ID 38956)
$wsALC.Range("H17").Value = 2255.2856
$wsALC.Range("N17").Value = -7607.000100000001
$wsALC.Range("L17").Value = 7271.000100000001
$wsALC.Range("J17").Value = 2423.6667
# Row 38 (Leve Item ID 4599)
$wsALC.Range("H38").Value = 16008.286
$wsALC.Range("J38").Value = 3800
$wsALC.Range("N38").Value = -12144
$wsALC.Range("L38").Value = 11400
$wsALC.Range("K38").Value = 75493.5
$wsALC.Range("M38").Value = -75121.5
$wsALC.Range("I38").Value = 25164.5
# Row 43 (Leve Item ID 5472)
$wsALC.Range("M43").Value = -1929.5
$wsALC.Range("H43").Value = 1998.75
$wsALC.Range("I43").Value = 1998.5
$wsALC.Range("L43").Value = 1999
$wsALC.Range("J43").Value = 1999
$wsALC.Range("N43").Value = -2137
$wsALC.Range("K43").Value = 1998.5
# Row 58 (Leve Item ID 4606)
$wsALC.Range("I58").Value = 1743.4
$wsALC.Range("H58").Value = 2802.4285
$wsALC.Range("L58").Value = 16350
$wsALC.Range("J58").Value = 5450
$wsALC.Range("N58").Value = -16650
$wsALC.Range("K58").Value = 5230.200000000001
$wsALC.Range("M58").Value = -5080.200000000001
# Row 62 (Leve Item ID 27781)
$wsALC.Range("I62").Value = 2847.5
$wsALC.Range("J62").Value = 2430
$wsALC.Range("H62").Value = 2708.3333
$wsALC.Range("N62").Value = -3678
$wsALC.Range("L62").Value = 2430
$wsALC.Range("K62").Value = 2847.5
$wsALC.Range("M62").Value = -2223.5
# Row 65 (Leve Item ID 27781)
$wsALC.Range("H65").Value = 2708.3333
$wsALC.Range("N65").Value = -18390
$wsALC.Range("L65").Value = 12150
$wsALC.Range("J65").Value = 2430
$wsALC.Range("K65").Value = 14237.5
$wsALC.Range("M65").Value = -11117.5
$wsALC.Range("I65").Value = 2847.5
# Row 100 (Leve Item ID 19906)
$wsALC.Range("H100").Value = 2214.6667
$wsALC.Range("I100").Value = 2214.6667
$wsALC.Range("K100").Value = 2214.6667
$wsALC.Range("M100").Value = -1673.6667
# Row 138 (Leve Item ID 44169)
$wsALC.Range("M138").Value = 38.9997999999996
$wsALC.Range("I138").Value = 1700.3334
$wsALC.Range("H138").Value = 2659.6726
$wsALC.Range("K138").Value = 5101.0002

$wsARM = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$wsARM.Range("K32").Value = 2783.0435
$wsARM.Range("M32").Value = -2496.0435
$wsARM.Range("H32").Value = 3593
$wsARM.Range("I32").Value = 2783.0435
$wsARM.Range("L32").Value = 22222
$wsARM.Range("J32").Value = 22222
$wsARM.Range("N32").Value = -22796
# Row 43 (Leve Item ID 21715)
$wsARM.Range("H43").Value = 38999
$wsARM.Range("L43").Value = 38999
$wsARM.Range("J43").Value = 38999
$wsARM.Range("N43").Value = -39625
# Row 45 (Leve Item ID 27714)
$wsARM.Range("I45").Value = 43512.6
$wsARM.Range("H45").Value = 27961.25
$wsARM.Range("N45").Value = -2796.3334
$wsARM.Range("L45").Value = 2042.3334
$wsARM.Range("J45").Value = 2042.3334
$wsARM.Range("K45").Value = 43512.6
$wsARM.Range("M45").Value = -43135.6
# Row 74 (Leve Item ID 44000)
$wsARM.Range("H74").Value = 3085.5
$wsARM.Range("J74").Value = 4245.9
$wsARM.Range("N74").Value = -5993.9
$wsARM.Range("L74").Value = 4245.9
# Row 77 (Leve Item ID 44000)
$wsARM.Range("L77").Value = 21229.5
$wsARM.Range("H77").Value = 3085.5
$wsARM.Range("J77").Value = 4245.9
$wsARM.Range("N77").Value = -29965.5
# Row 97 (Leve Item ID 19941)
$wsARM.Range("I97").Value = 1206.3334
$wsARM.Range("H97").Value = 1136.6154
$wsARM.Range("K97").Value = 1206.3334
$wsARM.Range("M97").Value = -710.3334
# Row 102 (Leve Item ID 19945)
$wsARM.Range("M102").Value = -1109.4443
$wsARM.Range("I102").Value = 2731.4443
$wsARM.Range("H102").Value = 3642.625
$wsARM.Range("K102").Value = 2731.4443
# Row 110 (Leve Item ID 27708)
$wsARM.Range("I110").Value = 3037.2
$wsARM.Range("H110").Value = 2847.28
$wsARM.Range("L110").Value = 2087.6
$wsARM.Range("J110").Value = 2087.6
$wsARM.Range("N110").Value = -6177.6
$wsARM.Range("K110").Value = 3037.2
$wsARM.Range("M110").Value = -992.1999999999998

$wsBSM = $wb.Worksheets.Item("BSM")
# Row 94 (Leve Item ID 19939)
$wsBSM.Range("H94").Value = 68966696
$wsBSM.Range("J94").Value = 2800
$wsBSM.Range("N94").Value = -3702
$wsBSM.Range("L94").Value = 2800
$wsBSM.Range("K94").Value = 76924060
$wsBSM.Range("M94").Value = -76923609
$wsBSM.Range("I94").Value = 76924060
# Row 105 (Leve Item ID 19947)
$wsBSM.Range("L105").Value = 22729686
$wsBSM.Range("J105").Value = 22729686
$wsBSM.Range("N105").Value = -22733180
$wsBSM.Range("H105").Value = 13686425
# Row 107 (Leve Item ID 27706)
$wsBSM.Range("M107").Value = -10988277
$wsBSM.Range("H107").Value = 5918338
$wsBSM.Range("I107").Value = 10990197
$wsBSM.Range("K107").Value = 10990197

$wsCRP = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$wsCRP.Range("I31").Value = 2132.3572
$wsCRP.Range("L31").Value = 5868.923
$wsCRP.Range("H31").Value = 3317.122
$wsCRP.Range("J31").Value = 5868.923
$wsCRP.Range("N31").Value = -6458.923
$wsCRP.Range("K31").Value = 2132.3572
$wsCRP.Range("M31").Value = -1837.3572
# Row 34 (Leve Item ID 44023)
$wsCRP.Range("M34").Value = -1930.3572
$wsCRP.Range("I34").Value = 2132.3572
$wsCRP.Range("L34").Value = 5868.923
$wsCRP.Range("H34").Value = 3317.122
$wsCRP.Range("J34").Value = 5868.923
$wsCRP.Range("N34").Value = -6272.923
$wsCRP.Range("K34").Value = 2132.3572
# Row 99 (Leve Item ID 36198)
$wsCRP.Range("H99").Value = 4375
$wsCRP.Range("K99").Value = 2000
$wsCRP.Range("M99").Value = -502
$wsCRP.Range("I99").Value = 2000
# Row 103 (Leve Item ID 19558)
$wsCRP.Range("K103").Value = 11172.333
$wsCRP.Range("H103").Value = 11172.333
$wsCRP.Range("M103").Value = -10000.333
$wsCRP.Range("I103").Value = 11172.333
# Row 126 (Leve Item ID 36198)
$wsCRP.Range("H126").Value = 4375
$wsCRP.Range("K126").Value = 6000
$wsCRP.Range("M126").Value = -3530
$wsCRP.Range("I126").Value = 2000
# Row 132 (Leve Item ID 44019)
$wsCRP.Range("K132").Value = 8324.6469
$wsCRP.Range("I132").Value = 2774.8823
$wsCRP.Range("M132").Value = -5794.6469
$wsCRP.Range("H132").Value = 11908943
$wsCRP.Range("L132").Value = 90928152
$wsCRP.Range("J132").Value = 30309384
$wsCRP.Range("N132").Value = -90933212
# Row 134 (Leve Item ID 44020)
$wsCRP.Range("I134").Value = 2612.8
$wsCRP.Range("H134").Value = 2681.303
$wsCRP.Range("J134").Value = 2895.375
$wsCRP.Range("N134").Value = -13756.125
$wsCRP.Range("L134").Value = 8686.125
$wsCRP.Range("K134").Value = 7838.400000000001
$wsCRP.Range("M134").Value = -5303.400000000001

$wsCUL = $wb.Worksheets.Item("CUL")
# Row 13 (Leve Item ID 4657)
$wsCUL.Range("L13").Value = 570
$wsCUL.Range("I13").Value = 221.25
$wsCUL.Range("H13").Value = 215
$wsCUL.Range("J13").Value = 190
$wsCUL.Range("N13").Value = -906
$wsCUL.Range("K13").Value = 663.75
$wsCUL.Range("M13").Value = -495.75
# Row 39 (Leve Item ID 4712)
$wsCUL.Range("I39").Value = 868.5
$wsCUL.Range("H39").Value = 2733.75
$wsCUL.Range("N39").Value = -14385
$wsCUL.Range("J39").Value = 4599
$wsCUL.Range("L39").Value = 13797
$wsCUL.Range("K39").Value = 2605.5
$wsCUL.Range("M39").Value = -2311.5

$wsGSM = $wb.Worksheets.Item("GSM")
# Row 7 (Leve Item ID 4197)
$wsGSM.Range("K7").Value = 11002
$wsGSM.Range("M7").Value = -10890
$wsGSM.Range("I7").Value = 11002
$wsGSM.Range("H7").Value = 11002
# Row 8 (Leve Item ID 4197)
$wsGSM.Range("K8").Value = 11002
$wsGSM.Range("M8").Value = -10863
$wsGSM.Range("I8").Value = 11002
$wsGSM.Range("H8").Value = 11002
# Row 107 (Leve Item ID 27802)
$wsGSM.Range("M107").Value = 1519.33334
$wsGSM.Range("H107").Value = 9067.375
$wsGSM.Range("I107").Value = 400.66666
$wsGSM.Range("K107").Value = 400.66666
# Row 141 (Leve Item ID 42504)
$wsGSM.Range("H141").Value = 52900
$wsGSM.Range("L141").Value = 52900
$wsGSM.Range("J141").Value = 52900
$wsGSM.Range("N141").Value = -63260

$wsLTW = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$wsLTW.Range("I46").Value = 2690.4285
$wsLTW.Range("J46").Value = 5844.5
$wsLTW.Range("N46").Value = -6220.5
$wsLTW.Range("H46").Value = 3636.65
$wsLTW.Range("L46").Value = 5844.5
$wsLTW.Range("K46").Value = 2690.4285
$wsLTW.Range("M46").Value = -2502.4285
# Row 64 (Leve Item ID 10810)
$wsLTW.Range("H64").Value = 39332.832
$wsLTW.Range("J64").Value = 39332.832
$wsLTW.Range("N64").Value = -39782.832
$wsLTW.Range("L64").Value = 39332.832
# Row 67 (Leve Item ID 10810)
$wsLTW.Range("J67").Value = 39332.832
$wsLTW.Range("N67").Value = -40892.832
$wsLTW.Range("H67").Value = 39332.832
$wsLTW.Range("L67").Value = 39332.832
# Row 136 (Leve Item ID 44060)
$wsLTW.Range("M136").Value = -5349
$wsLTW.Range("I136").Value = 2633
$wsLTW.Range("H136").Value = 2767.4736
$wsLTW.Range("K136").Value = 7899

$wsWVR = $wb.Worksheets.Item("WVR")
# Row 17 (Leve Item ID 3539)
$wsWVR.Range("I17").Value = 3874.25
$wsWVR.Range("H17").Value = 5099.4
$wsWVR.Range("K17").Value = 3874.25
$wsWVR.Range("M17").Value = -3702.25
# Row 97 (Leve Item ID 18220)
$wsWVR.Range("I97").Value = 0
$wsWVR.Range("H97").Value = 40000
$wsWVR.Range("K97").Value = 0
$wsWVR.Range("M97").ClearContents()
# Row 98 (Leve Item ID 18374)
$wsWVR.Range("H98").Value = 0
$wsWVR.Range("J98").Value = 0
$wsWVR.Range("N98").ClearContents()
$wsWVR.Range("L98").Value = 0
# Row 131 (Leve Item ID 34723)
$wsWVR.Range("J131").Value = 0
$wsWVR.Range("H131").Value = 0
$wsWVR.Range("N131").ClearContents()
$wsWVR.Range("L131").Value = 0
# Row 136 (Leve Item ID 44031)
$wsWVR.Range("M136").Value = -6367.200000000001
$wsWVR.Range("I136").Value = 2972.4
$wsWVR.Range("H136").Value = 3517.2183
$wsWVR.Range("N136").Value = -20010.201
$wsWVR.Range("L136").Value = 14910.201
$wsWVR.Range("J136").Value = 4970.067
$wsWVR.Range("K136").Value = 8917.200000000001
